$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the audit data (rows 2-254) with the GenEd dataset, growing the used range from A1:B187 to A1:B254.
$arr = New-Object 'object[,]' 253,2
$arr[0,0] = "GenEd---First Year Writing"
$arr[0,1] = "76-101"
$arr[1,0] = "GenEd---First Year Writing"
$arr[1,1] = "76-102"
$arr[2,0] = "GenEd---First Year Writing---2 Writing Minis"
$arr[2,1] = "76-106"
$arr[3,0] = "GenEd---First Year Writing---2 Writing Minis"
$arr[3,1] = "76-107"
$arr[4,0] = "GenEd---First Year Writing---2 Writing Minis"
$arr[4,1] = "76-108"
$arr[5,0] = "GenEd---Category 1---Category 1: Cognition, Choice, and Behavior (CS, CB, & HCI)"
$arr[5,1] = "70-311"
$arr[6,0] = "GenEd---Category 1---Category 1: Cognition, Choice, and Behavior (CS, CB, & HCI)"
$arr[6,1] = "70-318"
$arr[7,0] = "GenEd---Category 1---Category 1: Cognition, Choice, and Behavior (CS, CB, & HCI)"
$arr[7,1] = "70-385"
$arr[8,0] = "GenEd---Category 1---Category 1: Cognition, Choice, and Behavior (CS, CB, & HCI)"
$arr[8,1] = "80-101"
$arr[9,0] = "GenEd---Category 1---Category 1: Cognition, Choice, and Behavior (CS, CB, & HCI)"
$arr[9,1] = "80-130"
$arr[10,0] = "GenEd---Category 1---Category 1: Cognition, Choice, and Behavior (CS, CB, & HCI)"
$arr[10,1] = "80-150"
$arr[11,0] = "GenEd---Category 1---Category 1: Cognition, Choice, and Behavior (CS, CB, & HCI)"
$arr[11,1] = "80-180"
$arr[12,0] = "GenEd---Category 1---Category 1: Cognition, Choice, and Behavior (CS, CB, & HCI)"
$arr[12,1] = "80-221"
$arr[13,0] = "GenEd---Category 1---Category 1: Cognition, Choice, and Behavior (CS, CB, & HCI)"
$arr[13,1] = "80-252"
$arr[14,0] = "GenEd---Category 1---Category 1: Cognition, Choice, and Behavior (CS, CB, & HCI)"
$arr[14,1] = "80-270"
$arr[15,0] = "GenEd---Category 1---Category 1: Cognition, Choice, and Behavior (CS, CB, & HCI)"
$arr[15,1] = "80-271"
$arr[16,0] = "GenEd---Category 1---Category 1: Cognition, Choice, and Behavior (CS, CB, & HCI)"
$arr[16,1] = "80-275"
$arr[17,0] = "GenEd---Category 1---Category 1: Cognition, Choice, and Behavior (CS, CB, & HCI)"
$arr[17,1] = "80-330"
$arr[18,0] = "GenEd---Category 1---Category 1: Cognition, Choice, and Behavior (CS, CB, & HCI)"
$arr[18,1] = "85-102"
$arr[19,0] = "GenEd---Category 1---Category 1: Cognition, Choice, and Behavior (CS, CB, & HCI)"
$arr[19,1] = "85-104"
$arr[20,0] = "GenEd---Category 1---Category 1: Cognition, Choice, and Behavior (CS, CB, & HCI)"
$arr[20,1] = "85-211"
$arr[21,0] = "GenEd---Category 1---Category 1: Cognition, Choice, and Behavior (CS, CB, & HCI)"
$arr[21,1] = "85-213"
$arr[22,0] = "GenEd---Category 1---Category 1: Cognition, Choice, and Behavior (CS, CB, & HCI)"
$arr[22,1] = "85-221"
$arr[23,0] = "GenEd---Category 1---Category 1: Cognition, Choice, and Behavior (CS, CB, & HCI)"
$arr[23,1] = "85-241"
$arr[24,0] = "GenEd---Category 1---Category 1: Cognition, Choice, and Behavior (CS, CB, & HCI)"
$arr[24,1] = "85-251"
$arr[25,0] = "GenEd---Category 1---Category 1: Cognition, Choice, and Behavior (CS, CB, & HCI)"
$arr[25,1] = "85-261"
$arr[26,0] = "GenEd---Category 1---Category 1: Cognition, Choice, and Behavior (CS, CB, & HCI)"
$arr[26,1] = "85-370"
$arr[27,0] = "GenEd---Category 1---Category 1: Cognition, Choice, and Behavior (CS, CB, & HCI)"
$arr[27,1] = "85-408"
$arr[28,0] = "GenEd---Category 1---Category 1: Cognition, Choice, and Behavior (CS, CB, & HCI)"
$arr[28,1] = "85-414"
$arr[29,0] = "GenEd---Category 1---Category 1: Cognition, Choice, and Behavior (CS, CB, & HCI)"
$arr[29,1] = "85-421"
$arr[30,0] = "GenEd---Category 1---Category 1: Cognition, Choice, and Behavior (CS, CB, & HCI)"
$arr[30,1] = "88-120"
$arr[31,0] = "GenEd---Category 1---Category 1: Cognition, Choice, and Behavior (CS, CB, & HCI)"
$arr[31,1] = "88-230"
$arr[32,0] = "GenEd---Category 1---Category 1: Cognition, Choice, and Behavior (CS, CB, & HCI)"
$arr[32,1] = "88-231"
$arr[33,0] = "GenEd---Category 1---Category 1A: Cognitive Studies (AI)"
$arr[33,1] = "85-211"
$arr[34,0] = "GenEd---Category 1---Category 1A: Cognitive Studies (AI)"
$arr[34,1] = "85-213"
$arr[35,0] = "GenEd---Category 1---Category 1A: Cognitive Studies (AI)"
$arr[35,1] = "85-370"
$arr[36,0] = "GenEd---Category 1---Category 1A: Cognitive Studies (AI)"
$arr[36,1] = "85-408"
$arr[37,0] = "GenEd---Category 1---Category 1A: Cognitive Studies (AI)"
$arr[37,1] = "85-421"
$arr[38,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[38,1] = "19-101"
$arr[39,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[39,1] = "36-303"
$arr[40,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[40,1] = "66-221"
$arr[41,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[41,1] = "70-332"
$arr[42,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[42,1] = "73-102"
$arr[43,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[43,1] = "73-103"
$arr[44,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[44,1] = "73-104"
$arr[45,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[45,1] = "73-230"
$arr[46,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[46,1] = "73-240"
$arr[47,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[47,1] = "73-369"
$arr[48,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[48,1] = "76-425"
$arr[49,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[49,1] = "79-101"
$arr[50,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[50,1] = "79-155"
$arr[51,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[51,1] = "79-189"
$arr[52,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[52,1] = "79-212"
$arr[53,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[53,1] = "79-237"
$arr[54,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[54,1] = "79-244"
$arr[55,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[55,1] = "79-253"
$arr[56,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[56,1] = "79-275"
$arr[57,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[57,1] = "79-279"
$arr[58,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[58,1] = "79-300"
$arr[59,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[59,1] = "79-310"
$arr[60,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[60,1] = "79-315"
$arr[61,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[61,1] = "79-320"
$arr[62,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[62,1] = "79-321"
$arr[63,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[63,1] = "79-331"
$arr[64,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[64,1] = "79-343"
$arr[65,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[65,1] = "79-370"
$arr[66,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[66,1] = "79-383"
$arr[67,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[67,1] = "79-391"
$arr[68,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[68,1] = "79-392"
$arr[69,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[69,1] = "80-135"
$arr[70,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[70,1] = "80-136"
$arr[71,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[71,1] = "80-244"
$arr[72,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[72,1] = "80-245"
$arr[73,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[73,1] = "80-324"
$arr[74,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[74,1] = "80-334"
$arr[75,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[75,1] = "80-335"
$arr[76,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[76,1] = "80-348"
$arr[77,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[77,1] = "82-208"
$arr[78,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[78,1] = "84-104"
$arr[79,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[79,1] = "84-110"
$arr[80,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[80,1] = "84-275"
$arr[81,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[81,1] = "84-310"
$arr[82,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[82,1] = "84-322"
$arr[83,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[83,1] = "84-324"
$arr[84,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[84,1] = "84-352"
$arr[85,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[85,1] = "84-362"
$arr[86,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[86,1] = "84-365"
$arr[87,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[87,1] = "84-380"
$arr[88,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[88,1] = "84-386"
$arr[89,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[89,1] = "84-387"
$arr[90,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[90,1] = "84-389"
$arr[91,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[91,1] = "84-390"
$arr[92,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[92,1] = "84-393"
$arr[93,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[93,1] = "84-402"
$arr[94,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[94,1] = "84-405"
$arr[95,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[95,1] = "88-234"
$arr[96,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[96,1] = "88-255"
$arr[97,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[97,1] = "88-281"
$arr[98,0] = "GenEd---Category 2: Economic, Political, and Social Institutions"
$arr[98,1] = "88-284"
$arr[99,0] = "GenEd---Category 3: Cultural Analysis"
$arr[99,1] = "48-240"
$arr[100,0] = "GenEd---Category 3: Cultural Analysis"
$arr[100,1] = "48-241"
$arr[101,0] = "GenEd---Category 3: Cultural Analysis"
$arr[101,1] = "57-173"
$arr[102,0] = "GenEd---Category 3: Cultural Analysis"
$arr[102,1] = "60-105"
$arr[103,0] = "GenEd---Category 3: Cultural Analysis"
$arr[103,1] = "60-106"
$arr[104,0] = "GenEd---Category 3: Cultural Analysis"
$arr[104,1] = "62-371"
$arr[105,0] = "GenEd---Category 3: Cultural Analysis"
$arr[105,1] = "70-342"
$arr[106,0] = "GenEd---Category 3: Cultural Analysis"
$arr[106,1] = "76-221"
$arr[107,0] = "GenEd---Category 3: Cultural Analysis"
$arr[107,1] = "76-230"
$arr[108,0] = "GenEd---Category 3: Cultural Analysis"
$arr[108,1] = "76-232"
$arr[109,0] = "GenEd---Category 3: Cultural Analysis"
$arr[109,1] = "76-239"
$arr[110,0] = "GenEd---Category 3: Cultural Analysis"
$arr[110,1] = "76-241"
$arr[111,0] = "GenEd---Category 3: Cultural Analysis"
$arr[111,1] = "76-243"
$arr[112,0] = "GenEd---Category 3: Cultural Analysis"
$arr[112,1] = "76-339"
$arr[113,0] = "GenEd---Category 3: Cultural Analysis"
$arr[113,1] = "76-386"
$arr[114,0] = "GenEd---Category 3: Cultural Analysis"
$arr[114,1] = "79-104"
$arr[115,0] = "GenEd---Category 3: Cultural Analysis"
$arr[115,1] = "79-145"
$arr[116,0] = "GenEd---Category 3: Cultural Analysis"
$arr[116,1] = "79-201"
$arr[117,0] = "GenEd---Category 3: Cultural Analysis"
$arr[117,1] = "79-202"
$arr[118,0] = "GenEd---Category 3: Cultural Analysis"
$arr[118,1] = "79-211"
$arr[119,0] = "GenEd---Category 3: Cultural Analysis"
$arr[119,1] = "79-223"
$arr[120,0] = "GenEd---Category 3: Cultural Analysis"
$arr[120,1] = "79-226"
$arr[121,0] = "GenEd---Category 3: Cultural Analysis"
$arr[121,1] = "79-229"
$arr[122,0] = "GenEd---Category 3: Cultural Analysis"
$arr[122,1] = "79-230"
$arr[123,0] = "GenEd---Category 3: Cultural Analysis"
$arr[123,1] = "79-234"
$arr[124,0] = "GenEd---Category 3: Cultural Analysis"
$arr[124,1] = "79-240"
$arr[125,0] = "GenEd---Category 3: Cultural Analysis"
$arr[125,1] = "79-241"
$arr[126,0] = "GenEd---Category 3: Cultural Analysis"
$arr[126,1] = "79-242"
$arr[127,0] = "GenEd---Category 3: Cultural Analysis"
$arr[127,1] = "79-245"
$arr[128,0] = "GenEd---Category 3: Cultural Analysis"
$arr[128,1] = "79-248"
$arr[129,0] = "GenEd---Category 3: Cultural Analysis"
$arr[129,1] = "79-261"
$arr[130,0] = "GenEd---Category 3: Cultural Analysis"
$arr[130,1] = "79-262"
$arr[131,0] = "GenEd---Category 3: Cultural Analysis"
$arr[131,1] = "79-265"
$arr[132,0] = "GenEd---Category 3: Cultural Analysis"
$arr[132,1] = "79-281"
$arr[133,0] = "GenEd---Category 3: Cultural Analysis"
$arr[133,1] = "79-282"
$arr[134,0] = "GenEd---Category 3: Cultural Analysis"
$arr[134,1] = "79-286"
$arr[135,0] = "GenEd---Category 3: Cultural Analysis"
$arr[135,1] = "79-288"
$arr[136,0] = "GenEd---Category 3: Cultural Analysis"
$arr[136,1] = "79-293"
$arr[137,0] = "GenEd---Category 3: Cultural Analysis"
$arr[137,1] = "79-316"
$arr[138,0] = "GenEd---Category 3: Cultural Analysis"
$arr[138,1] = "79-329"
$arr[139,0] = "GenEd---Category 3: Cultural Analysis"
$arr[139,1] = "79-345"
$arr[140,0] = "GenEd---Category 3: Cultural Analysis"
$arr[140,1] = "79-350"
$arr[141,0] = "GenEd---Category 3: Cultural Analysis"
$arr[141,1] = "79-378"
$arr[142,0] = "GenEd---Category 3: Cultural Analysis"
$arr[142,1] = "79-386"
$arr[143,0] = "GenEd---Category 3: Cultural Analysis"
$arr[143,1] = "79-393"
$arr[144,0] = "GenEd---Category 3: Cultural Analysis"
$arr[144,1] = "79-395"
$arr[145,0] = "GenEd---Category 3: Cultural Analysis"
$arr[145,1] = "79-396"
$arr[146,0] = "GenEd---Category 3: Cultural Analysis"
$arr[146,1] = "79-465"
$arr[147,0] = "GenEd---Category 3: Cultural Analysis"
$arr[147,1] = "80-100"
$arr[148,0] = "GenEd---Category 3: Cultural Analysis"
$arr[148,1] = "80-250"
$arr[149,0] = "GenEd---Category 3: Cultural Analysis"
$arr[149,1] = "80-251"
$arr[150,0] = "GenEd---Category 3: Cultural Analysis"
$arr[150,1] = "80-253"
$arr[151,0] = "GenEd---Category 3: Cultural Analysis"
$arr[151,1] = "80-254"
$arr[152,0] = "GenEd---Category 3: Cultural Analysis"
$arr[152,1] = "80-255"
$arr[153,0] = "GenEd---Category 3: Cultural Analysis"
$arr[153,1] = "80-261"
$arr[154,0] = "GenEd---Category 3: Cultural Analysis"
$arr[154,1] = "80-276"
$arr[155,0] = "GenEd---Category 3: Cultural Analysis"
$arr[155,1] = "82-119"
$arr[156,0] = "GenEd---Category 3: Cultural Analysis"
$arr[156,1] = "82-267"
$arr[157,0] = "GenEd---Category 3: Cultural Analysis"
$arr[157,1] = "82-273"
$arr[158,0] = "GenEd---Category 3: Cultural Analysis"
$arr[158,1] = "82-279"
$arr[159,0] = "GenEd---Category 3: Cultural Analysis"
$arr[159,1] = "82-280"
$arr[160,0] = "GenEd---Category 3: Cultural Analysis"
$arr[160,1] = "82-282"
$arr[161,0] = "GenEd---Category 3: Cultural Analysis"
$arr[161,1] = "82-283"
$arr[162,0] = "GenEd---Category 3: Cultural Analysis"
$arr[162,1] = "82-286"
$arr[163,0] = "GenEd---Category 3: Cultural Analysis"
$arr[163,1] = "82-293"
$arr[164,0] = "GenEd---Category 3: Cultural Analysis"
$arr[164,1] = "82-294"
$arr[165,0] = "GenEd---Category 3: Cultural Analysis"
$arr[165,1] = "82-303"
$arr[166,0] = "GenEd---Category 3: Cultural Analysis"
$arr[166,1] = "82-304"
$arr[167,0] = "GenEd---Category 3: Cultural Analysis"
$arr[167,1] = "82-313"
$arr[168,0] = "GenEd---Category 3: Cultural Analysis"
$arr[168,1] = "82-314"
$arr[169,0] = "GenEd---Category 3: Cultural Analysis"
$arr[169,1] = "82-327"
$arr[170,0] = "GenEd---Category 3: Cultural Analysis"
$arr[170,1] = "82-331"
$arr[171,0] = "GenEd---Category 3: Cultural Analysis"
$arr[171,1] = "82-333"
$arr[172,0] = "GenEd---Category 3: Cultural Analysis"
$arr[172,1] = "82-342"
$arr[173,0] = "GenEd---Category 3: Cultural Analysis"
$arr[173,1] = "82-343"
$arr[174,0] = "GenEd---Category 3: Cultural Analysis"
$arr[174,1] = "82-344"
$arr[175,0] = "GenEd---Category 3: Cultural Analysis"
$arr[175,1] = "82-345"
$arr[176,0] = "GenEd---Category 3: Cultural Analysis"
$arr[176,1] = "82-436"
$arr[177,0] = "GenEd---Humanities/Arts Electives"
$arr[177,1] = "07-135"
$arr[178,0] = "GenEd---Humanities/Arts Electives"
$arr[178,1] = "11-423"
$arr[179,0] = "GenEd---Humanities/Arts Electives"
$arr[179,1] = "16-161"
$arr[180,0] = "GenEd---Humanities/Arts Electives"
$arr[180,1] = "16-397"
$arr[181,0] = "GenEd---Humanities/Arts Electives"
$arr[181,1] = "17-333"
$arr[182,0] = "GenEd---Humanities/Arts Electives"
$arr[182,1] = "17-562"
$arr[183,0] = "GenEd---Humanities/Arts Electives"
$arr[183,1] = "19-101"
$arr[184,0] = "GenEd---Humanities/Arts Electives"
$arr[184,1] = "19-351"
$arr[185,0] = "GenEd---Humanities/Arts Electives"
$arr[185,1] = "19-402"
$arr[186,0] = "GenEd---Humanities/Arts Electives"
$arr[186,1] = "19-403"
$arr[187,0] = "GenEd---Humanities/Arts Electives"
$arr[187,1] = "19-411"
$arr[188,0] = "GenEd---Humanities/Arts Electives"
$arr[188,1] = "21-150"
$arr[189,0] = "GenEd---Humanities/Arts Electives"
$arr[189,1] = "32-201"
$arr[190,0] = "GenEd---Humanities/Arts Electives"
$arr[190,1] = "32-402"
$arr[191,0] = "GenEd---Humanities/Arts Electives"
$arr[191,1] = "36-303"
$arr[192,0] = "GenEd---Humanities/Arts Electives"
$arr[192,1] = "70-100"
$arr[193,0] = "GenEd---Humanities/Arts Electives"
$arr[193,1] = "70-311"
$arr[194,0] = "GenEd---Humanities/Arts Electives"
$arr[194,1] = "70-318"
$arr[195,0] = "GenEd---Humanities/Arts Electives"
$arr[195,1] = "70-321"
$arr[196,0] = "GenEd---Humanities/Arts Electives"
$arr[196,1] = "70-332"
$arr[197,0] = "GenEd---Humanities/Arts Electives"
$arr[197,1] = "70-340"
$arr[198,0] = "GenEd---Humanities/Arts Electives"
$arr[198,1] = "70-341"
$arr[199,0] = "GenEd---Humanities/Arts Electives"
$arr[199,1] = "70-342"
$arr[200,0] = "GenEd---Humanities/Arts Electives"
$arr[200,1] = "70-343"
$arr[201,0] = "GenEd---Humanities/Arts Electives"
$arr[201,1] = "70-345"
$arr[202,0] = "GenEd---Humanities/Arts Electives"
$arr[202,1] = "70-348"
$arr[203,0] = "GenEd---Humanities/Arts Electives"
$arr[203,1] = "70-350"
$arr[204,0] = "GenEd---Humanities/Arts Electives"
$arr[204,1] = "70-352"
$arr[205,0] = "GenEd---Humanities/Arts Electives"
$arr[205,1] = "70-364"
$arr[206,0] = "GenEd---Humanities/Arts Electives"
$arr[206,1] = "70-365"
$arr[207,0] = "GenEd---Humanities/Arts Electives"
$arr[207,1] = "70-381"
$arr[208,0] = "GenEd---Humanities/Arts Electives"
$arr[208,1] = "70-430"
$arr[209,0] = "GenEd---Humanities/Arts Electives"
$arr[209,1] = "73-011"
$arr[210,0] = "GenEd---Humanities/Arts Electives"
$arr[210,1] = "73-102"
$arr[211,0] = "GenEd---Humanities/Arts Electives"
$arr[211,1] = "73-103"
$arr[212,0] = "GenEd---Humanities/Arts Electives"
$arr[212,1] = "99-231"
$arr[213,0] = "GenEd---Humanities/Arts Electives"
$arr[213,1] = "99-238"
$arr[214,0] = "GenEd---Science and Engineering---Science and Engineering (CS, AI, & HCI)---Science/Engineering, Any Department (4 courses)"
$arr[214,1] = "02-223"
$arr[215,0] = "GenEd---Science and Engineering---Science and Engineering (CS, AI, & HCI)---Science/Engineering, Any Department (4 courses)"
$arr[215,1] = "02-261"
$arr[216,0] = "GenEd---Science and Engineering---Science and Engineering (CS, AI, & HCI)---Science/Engineering, Any Department (4 courses)"
$arr[216,1] = "85-219"
$arr[217,0] = "GenEd---Science and Engineering---Science and Engineering (CS, AI, & HCI)---Science/Engineering, Any Department (4 courses)"
$arr[217,1] = "85-310"
$arr[218,0] = "GenEd---Science and Engineering---Science and Engineering (CS, AI, & HCI)---Science/Engineering, Any Department (4 courses)"
$arr[218,1] = "85-314"
$arr[219,0] = "GenEd---Science and Engineering---Science and Engineering (CS, AI, & HCI)---Science/Engineering, Same Department (2 courses)---Option 2"
$arr[219,1] = "02-223"
$arr[220,0] = "GenEd---Science and Engineering---Science and Engineering (CS, AI, & HCI)---Science/Engineering, Same Department (2 courses)---Option 2"
$arr[220,1] = "02-261"
$arr[221,0] = "GenEd---Science and Engineering---Science and Engineering (CS, AI, & HCI)---Science/Engineering, Same Department (2 courses)---Option 2"
$arr[221,1] = "02-262"
$arr[222,0] = "GenEd---Science and Engineering---Science and Engineering (CS, AI, & HCI)---Science/Engineering, Same Department (2 courses)---Option 2"
$arr[222,1] = "85-219"
$arr[223,0] = "GenEd---Science and Engineering---Science and Engineering (CS, AI, & HCI)---Lab Requirement"
$arr[223,1] = "02-261"
$arr[224,0] = "GenEd---Science and Engineering---Science and Engineering (CS, AI, & HCI)---Lab Requirement"
$arr[224,1] = "02-262"
$arr[225,0] = "GenEd---Science and Engineering---Science and Engineering (CS, AI, & HCI)---Lab Requirement"
$arr[225,1] = "03-124"
$arr[226,0] = "GenEd---Science and Engineering---Science and Engineering (CS, AI, & HCI)---Lab Requirement"
$arr[226,1] = "03-151"
$arr[227,0] = "GenEd---Science and Engineering---Science and Engineering (CS, AI, & HCI)---Lab Requirement"
$arr[227,1] = "03-206"
$arr[228,0] = "GenEd---Science and Engineering---Science and Engineering (CS, AI, & HCI)---Lab Requirement"
$arr[228,1] = "03-351"
$arr[229,0] = "GenEd---Science and Engineering---Science and Engineering (CS, AI, & HCI)---Lab Requirement"
$arr[229,1] = "09-101"
$arr[230,0] = "GenEd---Science and Engineering---Science and Engineering (CS, AI, & HCI)---Lab Requirement"
$arr[230,1] = "09-221"
$arr[231,0] = "GenEd---Science and Engineering---Science and Engineering (CS, AI, & HCI)---Lab Requirement"
$arr[231,1] = "27-100"
$arr[232,0] = "GenEd---Science and Engineering---Science and Engineering (CS, AI, & HCI)---Lab Requirement"
$arr[232,1] = "33-104"
$arr[233,0] = "GenEd---Science and Engineering---Science and Engineering (CS, AI, & HCI)---Lab Requirement"
$arr[233,1] = "33-228"
$arr[234,0] = "GenEd---Science and Engineering---Science and Engineering (CS, AI, & HCI)---Lab Requirement"
$arr[234,1] = "42-203"
$arr[235,0] = "GenEd---Science and Engineering---Science and Engineering (CS, AI, & HCI)---Lab Requirement"
$arr[235,1] = "85-310"
$arr[236,0] = "GenEd---Science and Engineering---Science and Engineering (CS, AI, & HCI)---Lab Requirement"
$arr[236,1] = "85-314"
$arr[237,0] = "GenEd---Science and Engineering---Science and Engineering (CB)---Science/Engineering, Any Department (4 courses)---Physics"
$arr[237,1] = "33-121"
$arr[238,0] = "GenEd---Science and Engineering---Science and Engineering (CB)---Science/Engineering, Any Department (4 courses)---Physics"
$arr[238,1] = "33-141"
$arr[239,0] = "GenEd---Science and Engineering---Science and Engineering (CB)---Science/Engineering, Any Department (4 courses)---Chemistry"
$arr[239,1] = "09-105"
$arr[240,0] = "GenEd---Science and Engineering---Science and Engineering (CB)---Science/Engineering, Any Department (4 courses)---Chemistry"
$arr[240,1] = "09-107"
$arr[241,0] = "GenEd---Science and Engineering---Science and Engineering (CB)---Science/Engineering, Any Department (4 courses)---Modern Biology"
$arr[241,1] = "03-121"
$arr[242,0] = "GenEd---Science and Engineering---Science and Engineering (CB)---Science/Engineering, Any Department (4 courses)---Modern Biology"
$arr[242,1] = "03-151"
$arr[243,0] = "GenEd---Science and Engineering---Science and Engineering (CB)---Science/Engineering, Any Department (4 courses)---Molecular Biology"
$arr[243,1] = "02-261"
$arr[244,0] = "GenEd---Science and Engineering---Science and Engineering (CB)---Science/Engineering, Any Department (4 courses)---Molecular Biology"
$arr[244,1] = "03-343"
$arr[245,0] = "GenEd---Science and Engineering---Science and Engineering (CB)---Science/Engineering, Same Department (2 courses)---Modern Biology"
$arr[245,1] = "03-121"
$arr[246,0] = "GenEd---Science and Engineering---Science and Engineering (CB)---Science/Engineering, Same Department (2 courses)---Modern Biology"
$arr[246,1] = "03-151"
$arr[247,0] = "GenEd---Science and Engineering---Science and Engineering (CB)---Science/Engineering, Same Department (2 courses)---Molecular Biology"
$arr[247,1] = "02-261"
$arr[248,0] = "GenEd---Science and Engineering---Science and Engineering (CB)---Science/Engineering, Same Department (2 courses)---Molecular Biology"
$arr[248,1] = "02-262"
$arr[249,0] = "GenEd---Science and Engineering---Science and Engineering (CB)---Science/Engineering, Same Department (2 courses)---Molecular Biology"
$arr[249,1] = "03-343"
$arr[250,0] = "GenEd---Science and Engineering---Science and Engineering (CB)---Lab Requirement"
$arr[250,1] = "02-261"
$arr[251,0] = "GenEd---Science and Engineering---Science and Engineering (CB)---Lab Requirement"
$arr[251,1] = "02-262"
$arr[252,0] = "GenEd---Science and Engineering---Science and Engineering (CB)---Lab Requirement"
$arr[252,1] = "03-343"

$ws.Range("A2:B254").Value = $arr

Write-Output "Wrote $($arr.GetLength(0)) data rows to A2:B254"
